# Updated cryptos list values/links per upstream diff.
# D-column price cells are stored as text in the source sheet even when
# numeric-looking (e.g. "0.999", "541.95") -- coerce with NumberFormat
# '@' before assigning, then ClearFormats() to drop the temporary text
# format/quote-prefix so the cell style matches the original (no style).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "65.856.98"
$ws.Range("E2").Value = "  +3.26%  "
$ws.Range("D3").Value = "2.668.44"
$ws.Range("E3").Value = "  +1.86%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +2.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.49"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +5.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -0.11%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.38%  "
$ws.Range("D9").Value = "2.667.20"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("E10").Value = "  +9.37%  "
$ws.Range("E11").Value = "  +2.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.87"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.43%  "
$ws.Range("E13").Value = "  +1.74%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.94"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +7.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000196"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +15.39%  "
$ws.Range("D16").Value = "3.141.11"
$ws.Range("E16").Value = "  +1.64%  "
$ws.Range("D17").Value = "65.538.66"
$ws.Range("E17").Value = "  +2.90%  "
$ws.Range("D18").Value = "2.653.00"
$ws.Range("E18").Value = "  +1.38%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.73"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.89"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "361.57"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.46"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +6.54%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.38"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.71"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.12%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.56"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +4.38%  "
$ws.Range("E27").Value = "  +17.55%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.65"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.70%  "
$ws.Range("E29").Value = "  +2.97%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.20"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("E31").Value = "  +6.43%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "541.95"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.01%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.00"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.82"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.06%  "
$ws.Range("E35").Value = "  +2.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.38"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +4.05%  "
$ws.Range("E37").Value = "  +3.42%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.73"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +4.08%  "
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.02"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.28%  "
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "162.52"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.34%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.44"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "166.22"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -1.05%  "
$ws.Range("E45").Value = "  +2.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +9.66%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0615"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +5.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.14"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -0.04%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.662"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +4.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0265"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.28%  "
$ws.Range("E51").Value = "  +1.77%  "
